$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question text in B5: remove the quote marks around the "+" symbol
$ws.Range("B5").Value = 'What basic operation does it have as a + symbol?'

# Update the active selection to B6 (as recorded in the sheet view)
$ws.Range("B6").Select()
